# "env can run now, but some parameter need to be modified by model browser"
#
# 1. Fix the absolute path recorded by Excel for this workbook (it now lives
#    one directory deeper, under "config\cfg_crlc" instead of "cfg_crlc").
# 2. The B-column parameter label for rows 2-6 needs to change from "cvc1"
#    to "cvc4" (rows 7-8 keep their existing "parab" label).
# 3. Move the active selection to E4 (was D9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the recorded absolute path ---------------------------------
$wb.AbsolutePath = "C:\digestion\tiny_ssu_calc\config\cfg_crlc\"

# --- 2. Update the parameter label for rows 2-6 ---------------------------
$ws.Range("B2:B6").Value = "cvc4"

# --- 3. Move the selection -------------------------------------------------
$ws.Range("E4").Select()
